$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I3").Value = 184
